$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sandeep"
$ws.Range("B2").Value = "Sandeep1@gmail.com"
$ws.Range("C2").Value = 7817004567
$ws.Range("D2").Value = "Vadodara"
$ws.Range("F2").Value = "Sandeep@123"
$ws.Range("G2").Value = "Sandeep@1234"

$ws.Range("H2").Select()
